# Bank.xlsx update — "ไม่ระบุเวลา" (time not specified) replaces "-" in the
# Saturday/Sunday hours columns (K:L), header style for K1:L1 is normalized to
# match the rest of the header row, and the sheet's scroll/selection state is
# updated to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replace "-" with "ไม่ระบุเวลา" in the Saturday (K) / Sunday (L) columns ---
$newText = "ไม่ระบุเวลา"
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("K", "L")) {
        $cell = $ws.Range("$col$r")
        if ($cell.Value2 -eq "-") {
            $cell.Value2 = $newText
        }
    }
}

# --- 2. Normalize K1:L1 header style to match the other header cells (A1:J1, M1:P1) ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Update the active window's scroll position / selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 103
$ws.Range("K137:L156").Select() | Out-Null

Write-Output "done"
